$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = (Get-Date -Year 2022 -Month 2 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 600
$ws.Range("L2").Value = 650
$ws.Range("M2").Value = 625
$ws.Range("P2").Value = 625

$ws.Range("D3").Value = (Get-Date -Year 2022 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K3").Value = 700
$ws.Range("L3").Value = 800
$ws.Range("M3").Value = 750
$ws.Range("P3").Value = 750

$ws.Range("D4").Value = (Get-Date -Year 2022 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = 600
$ws.Range("P4").Value = 600

$ws.Range("D5").Value = (Get-Date -Year 2022 -Month 3 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 650
$ws.Range("L5").Value = 700
$ws.Range("M5").Value = 675
$ws.Range("P5").Value = 675

$ws.Range("D6").Value = (Get-Date -Year 2022 -Month 7 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 700
$ws.Range("L6").Value = 800
$ws.Range("M6").Value = 750
$ws.Range("P6").Value = 750

$ws.Range("D7").Value = (Get-Date -Year 2022 -Month 7 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 600
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = 600
$ws.Range("P7").Value = 600

$ws.Range("D8").Value = (Get-Date -Year 2022 -Month 8 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 750
$ws.Range("L8").Value = 850
$ws.Range("M8").Value = 800
$ws.Range("P8").Value = 800

$ws.Range("D9").Value = (Get-Date -Year 2022 -Month 8 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K9").Value = 650
$ws.Range("L9").Value = 650
$ws.Range("M9").Value = 650
$ws.Range("P9").Value = 650

$ws.Range("D10").Value = (Get-Date -Year 2022 -Month 7 -Day 19 -Hour 0 -Minute 0 -Second 0)

$ws.Range("D11").Value = (Get-Date -Year 2022 -Month 7 -Day 19 -Hour 0 -Minute 0 -Second 0)

$ws.Range("D12").Value = (Get-Date -Year 2022 -Month 8 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J12").Value = 240
$ws.Range("K12").Value = 750
$ws.Range("L12").Value = 850
$ws.Range("M12").Value = 800
$ws.Range("P12").Value = 800

$ws.Range("D13").Value = (Get-Date -Year 2022 -Month 8 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I13").Value = "Segunda"
$ws.Range("K13").Value = 650
$ws.Range("L13").Value = 650
$ws.Range("M13").Value = 650
$ws.Range("P13").Value = 650

$ws.Range("D14").Value = (Get-Date -Year 2022 -Month 9 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 300
$ws.Range("K14").Value = 750
$ws.Range("L14").Value = 850
$ws.Range("M14").Value = 800
$ws.Range("P14").Value = 800

$ws.Range("D15").Value = (Get-Date -Year 2022 -Month 7 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L15").Value = 750
$ws.Range("M15").Value = 725
$ws.Range("P15").Value = 725

$ws.Range("D16").Value = (Get-Date -Year 2022 -Month 10 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 700
$ws.Range("L16").Value = 800
$ws.Range("M16").Value = 750
$ws.Range("P16").Value = 750

$ws.Range("D17").Value = (Get-Date -Year 2022 -Month 10 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I17").Value = "Segunda"
$ws.Range("J17").Value = 150
$ws.Range("L17").Value = 600
$ws.Range("M17").Value = 600
$ws.Range("P17").Value = 600

$ws.Range("D18").Value = (Get-Date -Year 2022 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J18").Value = 300

$ws.Range("D19").Value = (Get-Date -Year 2022 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0)

$ws.Range("D20").Value = (Get-Date -Year 2022 -Month 8 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J20").Value = 240
$ws.Range("K20").Value = 750
$ws.Range("M20").Value = 775
$ws.Range("P20").Value = 775

$ws.Range("D21").Value = (Get-Date -Year 2022 -Month 8 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J21").Value = 250
$ws.Range("K21").Value = 650
$ws.Range("L21").Value = 650
$ws.Range("M21").Value = 650
$ws.Range("P21").Value = 650

$ws.Range("D22").Value = (Get-Date -Year 2022 -Month 2 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 650
$ws.Range("M22").Value = 625
$ws.Range("P22").Value = 625

$ws.Range("D23").Value = (Get-Date -Year 2022 -Month 9 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 800
$ws.Range("L23").Value = 900
$ws.Range("M23").Value = 850
$ws.Range("P23").Value = 850

$ws.Range("D24").Value = (Get-Date -Year 2022 -Month 8 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J24").Value = 160

$ws.Range("D25").Value = (Get-Date -Year 2022 -Month 8 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I25").Value = "Segunda"
$ws.Range("J25").Value = 120
$ws.Range("K25").Value = 650
$ws.Range("M25").Value = 650
$ws.Range("P25").Value = 650
